$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1212.97
$ws.Range("I112").Value = 371.42856
$ws.Range("J112").Value = 1276.3118
$ws.Range("K112").Value = 1114.28568
$ws.Range("L112").Value = 3828.9354
$ws.Range("M112").Value = -6.285679999999957
$ws.Range("N112").Value = -6044.9354

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 4099662.5
$ws.Range("J129").Value = 1288.4108
$ws.Range("L129").Value = 3865.2324
$ws.Range("N129").Value = -13865.2324

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 5408190.5
$ws.Range("I132").Value = 5716924
$ws.Range("J132").Value = 5350
$ws.Range("K132").Value = 17150772
$ws.Range("L132").Value = 16050
$ws.Range("M132").Value = -17148242
$ws.Range("N132").Value = -21110

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2866518.2
$ws.Range("I137").Value = 3343488.2
$ws.Range("J137").Value = 4698
$ws.Range("K137").Value = 10030464.6
$ws.Range("L137").Value = 14094
$ws.Range("M137").Value = -10027914.6
$ws.Range("N137").Value = -19194

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 441040.1
$ws.Range("I141").Value = 1717.8846
$ws.Range("J141").Value = 2072808.2
$ws.Range("K141").Value = 5153.6538
$ws.Range("L141").Value = 6218424.6
$ws.Range("M141").Value = 26.34619999999995
$ws.Range("N141").Value = -6228784.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2659.9412
$ws.Range("I61").Value = 1279.3572
$ws.Range("J61").Value = 3626.35
$ws.Range("K61").Value = 1279.3572
$ws.Range("L61").Value = 3626.35
$ws.Range("M61").Value = -1067.3572
$ws.Range("N61").Value = -4050.35

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 759.56525
$ws.Range("I74").Value = 725.9091
$ws.Range("J74").Value = 1500
$ws.Range("K74").Value = 725.9091
$ws.Range("L74").Value = 1500
$ws.Range("M74").Value = 148.0909
$ws.Range("N74").Value = -3248

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 759.56525
$ws.Range("I77").Value = 725.9091
$ws.Range("J77").Value = 1500
$ws.Range("K77").Value = 3629.5455
$ws.Range("L77").Value = 7500
$ws.Range("M77").Value = 738.4545000000003
$ws.Range("N77").Value = -16236

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2659.9412
$ws.Range("I136").Value = 1279.3572
$ws.Range("J136").Value = 3626.35
$ws.Range("K136").Value = 3838.0716
$ws.Range("L136").Value = 10879.05
$ws.Range("M136").Value = -1288.0716
$ws.Range("N136").Value = -15979.05

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H137").Value = 44400
$ws.Range("J137").Value = 44400
$ws.Range("L137").Value = 44400
$ws.Range("N137").Value = -54600

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H139").Value = 27228.125
$ws.Range("J139").Value = 27228.125
$ws.Range("L139").Value = 27228.125
$ws.Range("N139").Value = -37508.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2923.9167
$ws.Range("I134").Value = 2480.25
$ws.Range("J134").Value = 3811.25
$ws.Range("K134").Value = 7440.75
$ws.Range("L134").Value = 11433.75
$ws.Range("M134").Value = -4905.75
$ws.Range("N134").Value = -16503.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2441216.2
$ws.Range("J31").Value = 3166.6667
$ws.Range("L31").Value = 3166.6667
$ws.Range("N31").Value = -3756.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2441216.2
$ws.Range("J34").Value = 3166.6667
$ws.Range("L34").Value = 3166.6667
$ws.Range("N34").Value = -3570.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 20125
$ws.Range("J52").Value = 20125
$ws.Range("L52").Value = 20125
$ws.Range("N52").Value = -20713

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 12502984
$ws.Range("I58").Value = 1991.4166
$ws.Range("J58").Value = 31254474
$ws.Range("K58").Value = 1991.4166
$ws.Range("L58").Value = 31254474
$ws.Range("M58").Value = -1788.4166
$ws.Range("N58").Value = -31254880

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H81").Value = 37600
$ws.Range("J81").Value = 37600
$ws.Range("L81").Value = 37600
$ws.Range("N81").Value = -39596

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H84").Value = 37600
$ws.Range("J84").Value = 37600
$ws.Range("L84").Value = 112800
$ws.Range("N84").Value = -122784

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2462.6743
$ws.Range("I132").Value = 1960.4814
$ws.Range("J132").Value = 3310.125
$ws.Range("K132").Value = 5881.4442
$ws.Range("L132").Value = 9930.375
$ws.Range("M132").Value = -3351.4442
$ws.Range("N132").Value = -14990.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1984.2
$ws.Range("I134").Value = 974.8
$ws.Range("J134").Value = 2993.6
$ws.Range("K134").Value = 2924.4
$ws.Range("L134").Value = 8980.799999999999
$ws.Range("M134").Value = -389.3999999999996
$ws.Range("N134").Value = -14050.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 12502984
$ws.Range("I136").Value = 1991.4166
$ws.Range("J136").Value = 31254474
$ws.Range("K136").Value = 5974.2498
$ws.Range("L136").Value = 93763422
$ws.Range("M136").Value = -3424.2498
$ws.Range("N136").Value = -93768522

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 5520.1055
$ws.Range("I3").Value = 2542.6667
$ws.Range("J3").Value = 6078.375
$ws.Range("K3").Value = 7628.000100000001
$ws.Range("L3").Value = 18235.125
$ws.Range("M3").Value = -7516.000100000001
$ws.Range("N3").Value = -18459.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 2200
$ws.Range("I75").Value = 1000
$ws.Range("K75").Value = 3000
$ws.Range("M75").Value = -2002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H78").Value = 2200
$ws.Range("I78").Value = 1000
$ws.Range("K78").Value = 9000
$ws.Range("M78").Value = -4008

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2817.4
$ws.Range("I126").Value = 1971.1
$ws.Range("J126").Value = 4510
$ws.Range("K126").Value = 5913.299999999999
$ws.Range("L126").Value = 13530
$ws.Range("M126").Value = -3443.299999999999
$ws.Range("N126").Value = -18470

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3037.5938
$ws.Range("I132").Value = 1706.0588
$ws.Range("J132").Value = 4546.6665
$ws.Range("K132").Value = 5118.1764
$ws.Range("L132").Value = 13639.9995
$ws.Range("M132").Value = -2588.1764
$ws.Range("N132").Value = -18699.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2860038
$ws.Range("I136").Value = 6669820
$ws.Range("J136").Value = 2701.8
$ws.Range("K136").Value = 20009460
$ws.Range("L136").Value = 8105.400000000001
$ws.Range("M136").Value = -20006910
$ws.Range("N136").Value = -13205.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 200012.88
$ws.Range("I132").Value = 264669.84
$ws.Range("J132").Value = 11015.615
$ws.Range("K132").Value = 794009.52
$ws.Range("L132").Value = 33046.845
$ws.Range("M132").Value = -791479.52
$ws.Range("N132").Value = -38106.845

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1114.8788
$ws.Range("I136").Value = 630.7241
$ws.Range("J136").Value = 4625
$ws.Range("K136").Value = 1892.1723
$ws.Range("L136").Value = 13875
$ws.Range("M136").Value = 657.8276999999998
$ws.Range("N136").Value = -18975
